$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.607.44"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.471.91"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +1.14%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0857"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "2.852.57"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.12%  "
$ws.Range("D16").Value = "2.484.98"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "41.556.29"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.97"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").Value = "1.990.81"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").Value = "2.709.92"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.58%  "
